# Updated cryptos list - price and volume(1h) refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.371.89"
$ws.Range("E2").Value = "  +3.97%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.492.81"
$ws.Range("E3").Value = "  +3.31%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.15"
$ws.Range("E5").Value = "  +3.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.04"
$ws.Range("E6").Value = "  +4.22%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.603"
$ws.Range("E8").Value = "  +11.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.492.30"
$ws.Range("E9").Value = "  +3.23%  "

$ws.Range("E10").Value = "  -0.83%  "

$ws.Range("E11").Value = "  +3.89%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.448"
$ws.Range("E12").Value = "  +3.95%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.093.39"
$ws.Range("E13").Value = "  +3.21%  "

$ws.Range("E14").Value = "  +1.23%  "

$ws.Range("E15").Value = "  +3.70%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.74"
$ws.Range("E16").Value = "  +6.98%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.371.84"
$ws.Range("E17").Value = "  +3.74%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.508.72"
$ws.Range("E18").Value = "  +4.56%  "

$ws.Range("E19").Value = "  +4.26%  "

$ws.Range("E20").Value = "  +3.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "386.40"
$ws.Range("E21").Value = "  +2.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.33"
$ws.Range("E22").Value = "  +3.71%  "

$ws.Range("E23").Value = "  +5.34%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.24"
$ws.Range("E24").Value = "  +2.74%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.77%  "

$ws.Range("E26").Value = "  +6.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.01"
$ws.Range("E27").Value = "  +6.33%  "

$ws.Range("E28").Value = "  +2.50%  "

$ws.Range("E29").Value = "  +0.16%  "

$ws.Range("E30").Value = "  +13.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.25"
$ws.Range("E31").Value = "  +4.52%  "

$ws.Range("E32").Value = "  +4.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.67"
$ws.Range("E33").Value = "  +3.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.34"
$ws.Range("E34").Value = "  +9.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.59"
$ws.Range("E35").Value = "  +10.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.64"
$ws.Range("E36").Value = "  +3.16%  "

$ws.Range("E37").Value = "  +6.86%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.017.25"
$ws.Range("E38").Value = "  +4.27%  "

$ws.Range("E39").Value = "  +2.94%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "27.27"
$ws.Range("E40").Value = "  +2.27%  "

$ws.Range("E41").Value = "  +3.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.60"
$ws.Range("E42").Value = "  +6.95%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.54"
$ws.Range("E43").Value = "  +3.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.83"
$ws.Range("E44").Value = "  +4.36%  "

$ws.Range("E45").Value = "  +3.92%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.67"
$ws.Range("E46").Value = "  +11.50%  "

$ws.Range("E47").Value = "  +5.71%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "326.18"
$ws.Range("E48").Value = "  +14.73%  "

$ws.Range("E49").Value = "  +6.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.21"
$ws.Range("E50").Value = "  +3.68%  "

$ws.Range("E51").Value = "  +6.83%  "
